$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Formula = "=18+24"

$ws.Range("D2").Select()
